$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.253.20'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = '1.592.40'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.500'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.71%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.245'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0607'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.96'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.31%  '
$ws.Range("E11").Value = '  +0.40%  '
$ws.Range("D12").Value = '1.816.30'
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").Value = '1.590.56'
$ws.Range("E13").Value = '  -0.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.05%  '
$ws.Range("E15").Value = '  -2.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.76'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("D17").Value = '26.249.92'
$ws.Range("E17").Value = '  -0.39%  '
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.30%  '
$ws.Range("E20").Value = '  -1.60%  '
$ws.Range("E22").Value = '  -0.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.63'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  -1.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.112'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.26%  '
$ws.Range("E29").Value = '  -0.68%  '
$ws.Range("E30").Value = '  -2.40%  '
$ws.Range("E31").Value = '  +0.55%  '
$ws.Range("E32").Value = '  -0.57%  '
$ws.Range("D33").Value = '1.408.71'
$ws.Range("E33").Value = '  +5.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.97'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("E35").Value = '  -0.52%  '
$ws.Range("E36").Value = '  -1.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.574'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.48%  '
$ws.Range("E38").Value = '  -1.15%  '
$ws.Range("E39").Value = '  +0.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.76'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.70%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.961'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.88%  '
$ws.Range("E43").Value = '  +0.96%  '
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("D45").Value = '1.728.19'
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.86'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '87.16'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.81%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0104'
$ws.Range("E48").Value = '  -0.76%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.49'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.45%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0502'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.63%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0952'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.38%  '
